$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 49 (pushes existing rows 49:66 down to 50:67)
$ws.Rows.Item(49).Insert()

# Populate the new row 49 with the weekly price-report entry
$ws.Range("A49").Value = 10
$ws.Range("B49").Value = "Vega Modelo de Temuco"
$ws.Range("C49").Value = "La Araucanía"
$ws.Range("D49").Value = 44830
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = 300000000
$ws.Range("G49").Value = "Espárragos"
$ws.Range("H49").Value = "Verde"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 85
$ws.Range("K49").Value = 29000
$ws.Range("L49").Value = 29000
$ws.Range("M49").Value = 29000
$ws.Range("N49").Value = "$/caja 10 kilos"
$ws.Range("O49").Value = "Provincia de Quillota"
$ws.Range("P49").Value = 2900
$ws.Range("Q49").Value = 10
$ws.Range("R49").Value = "Hortaliza"
